$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.322.61"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "3.160.07"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.94"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.49"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "3.157.35"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("E10").Value = "  -6.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -7.91%  "
$ws.Range("E12").Value = "  -6.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -7.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.53"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -9.63%  "
$ws.Range("D15").Value = "3.685.87"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "64.317.08"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "3.163.68"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  -6.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.37"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -5.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.71"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.66"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -8.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.67"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -6.88%  "
$ws.Range("E29").Value = "  -8.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.71"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("E31").Value = "  -21.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.46%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.25"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -6.73%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.56"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -7.39%  "
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("E38").Value = "  -9.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.87"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -8.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -12.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0395"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -7.73%  "
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("E43").Value = "  -8.32%  "
$ws.Range("D44").Value = "2.842.65"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("E45").Value = "  -9.61%  "
$ws.Range("E46").Value = "  -9.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.45"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -8.04%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -7.25%  "
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.57"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.58%  "
